$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix D3 (contacto for H2): 968693641 -> 968693643
$ws.Range("D3").Value = 968693643

# Add new row 5: H4 hotel
$ws.Range("A5").Value = "H4"
$ws.Range("B5").Value = 638
$ws.Range("C5").Value = "Porto"
$ws.Range("D5").Value = 968693641
$ws.Range("E5").Value = 80

# Update selection to match post-edit state (next empty row)
$ws.Range("E6").Select()
